# Kazakhstan Premier League 2023 update
# - Swap the per-match data (columns F:V) between three pairs of adjacent
#   rows (the two fixtures of the same matchday had been recorded in the
#   wrong order upstream).
# - Append two newly-scraped fixtures as rows 158 and 159.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap F:V between matched row pairs --------------------------------
# (named parameters are unreliable in this host, so plain script-level
# loops/arrays are used instead of helper functions)

$swapPairs = @(@(66, 67), @(115, 116), @(127, 128))

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($col = 6; $col -le 22; $col++) {
        $cellA = $ws.Cells.Item($r1, $col)
        $cellB = $ws.Cells.Item($r2, $col)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

# --- 2) Append two new fixtures -------------------------------------------

# Copy the number formatting from the last existing row (157) so the new
# rows look the same (bold+bordered index column, datetime match-date
# column); everything else uses the plain default style.
$ws.Cells.Item(157, 1).Copy()
$ws.Cells.Item(158, 1).PasteSpecial(-4122)
$ws.Cells.Item(157, 5).Copy()
$ws.Cells.Item(158, 5).PasteSpecial(-4122)
$ws.Cells.Item(157, 1).Copy()
$ws.Cells.Item(159, 1).PasteSpecial(-4122)
$ws.Cells.Item(157, 5).Copy()
$ws.Cells.Item(159, 5).PasteSpecial(-4122)

$row158 = @(157, "kazakhstan", "premier-league", "2023", 45196.625, "FC Astana", 5, "Okzhetpes", 2, 1.1, "04/08/2023 04:12", 1.18, "27/09/2023 14:55", 6.96, "04/08/2023 04:12", 5.89, "27/09/2023 14:57", 11.37, "04/08/2023 04:12", 11.66, "27/09/2023 14:57", "https://www.betexplorer.com/football/kazakhstan/premier-league/fc-astana-okzhetpes/QwANKdX5/")
$row159 = @(158, "kazakhstan", "premier-league", "2023", 45196.66666666666, "Ordabasy", 1, "Maqtaaral", 0, 1.23, "14/08/2023 15:21", 1.23, "14/08/2023 15:21", 5.28, "14/08/2023 15:21", 5.28, "14/08/2023 15:21", 9.09, "14/08/2023 15:21", 9.09, "14/08/2023 15:21", "https://www.betexplorer.com/football/kazakhstan/premier-league/ordabasy-maqtaaral/4U8RJGnC/")
$newRows = @($row158, $row159)

$destRows = @(158, 159)

for ($i = 0; $i -lt 2; $i++) {
    $destRow = $destRows[$i]
    $data = $newRows[$i]
    for ($col = 1; $col -le 22; $col++) {
        $ws.Cells.Item($destRow, $col).Value = $data[$col - 1]
    }
    # Column D ("temporada") is textual ("2023"); without the leading
    # apostrophe Excel's smart-typing would store it as the number 2023.
    # Style = "Normal" removes the resulting quote-prefix style so the
    # cell ends up with the same default formatting as its neighbours.
    $ws.Cells.Item($destRow, 4).Value = "'" + $data[3]
    $ws.Cells.Item($destRow, 4).Style = "Normal"
}

Write-Output "done"
